# RoutingPlan.xlsx edit — "login and register page done".
#
# Rows 4-8 of the plan (the "/", "/login" and "/register" routes) already
# have a "v" (done) mark in column F ("routers"); mirror that same mark
# into column G ("ejs page") now that those pages are implemented too.
# Then move the view/selection on to where work continues next.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

foreach ($row in 4..8) {
    $ws.Cells.Item($row, 7).Value = "v"
}

# Scroll the sheet down and land the selection on the next open row.
$excel.Goto($ws.Range("A7"), $true)
$ws.Range("E24").Select()
